# Update "想去人数" (interest count) values in the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 57
$ws1.Range("F4").Value = 156
$ws1.Range("F6").Value = 5152
$ws1.Range("F7").Value = 113
$ws1.Range("F8").Value = 5298
$ws1.Range("F9").Value = 611
$ws1.Range("F10").Value = 1346

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 57
$ws4.Range("F4").Value = 156
$ws4.Range("F7").Value = 5152
$ws4.Range("F8").Value = 113
$ws4.Range("F9").Value = 5298
$ws4.Range("F10").Value = 611
$ws4.Range("F11").Value = 1346
